$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.862.76'
$ws.Range("E2").Value = '  -0.73%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.938.34'
$ws.Range("E3").Value = '  -0.86%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.85'
$ws.Range("E5").Value = '  -0.86%  '

# Row 6
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4908'
$ws.Range("E7").Value = '  +0.27%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2944'
$ws.Range("E8").Value = '  -0.93%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06899'
$ws.Range("E9").Value = '  +0.73%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.27'
$ws.Range("E10").Value = '  +0.59%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '105.21'
$ws.Range("E11").Value = '  -2.63%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.948.02'
$ws.Range("E12").Value = '  +1.42%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07787'
$ws.Range("E13").Value = '  +0.39%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.363'
$ws.Range("E14").Value = '  -1.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7007'
$ws.Range("E15").Value = '  -1.04%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '275.09'
$ws.Range("E16").Value = '  -3.31%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.826.50'
$ws.Range("E17").Value = '  -0.89%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007733'
$ws.Range("E18").Value = '  -0.39%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.09'
$ws.Range("E19").Value = '  -1.04%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.596'
$ws.Range("E20").Value = '  +2.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").Value = '  -0.10%  '

# Row 22
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.171.50'
$ws.Range("E22").Value = '  -0.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9992'
$ws.Range("E23").Value = '  -0.19%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.526'
$ws.Range("E24").Value = '  +0.08%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.874'
$ws.Range("E25").Value = '  +0.71%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.21'
$ws.Range("E26").Value = '  -1.57%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.64'
$ws.Range("E27").Value = '  -2.25%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.165'
$ws.Range("E28").Value = '  -2.47%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1045'
$ws.Range("E29").Value = '  -1.07%  '

# Row 30
$ws.Range("E30").Value = '  -2.76%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.558'
$ws.Range("E31").Value = '  -1.80%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.574'
$ws.Range("E32").Value = '  -0.22%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.385'
$ws.Range("E33").Value = '  -1.35%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04887'
$ws.Range("E34").Value = '  -1.33%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7618'
$ws.Range("E35").Value = '  +1.01%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.152'
$ws.Range("E36").Value = '  -2.38%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9997'
$ws.Range("E37").Value = '  -0.08%  '

# Row 38
$ws.Range("E38").Value = '  +0.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02013'
$ws.Range("E39").Value = '  -1.15%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.655'
$ws.Range("E40").Value = '  -1.89%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.520'
$ws.Range("E41").Value = '  +1.35%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '78.36'
$ws.Range("E42").Value = '  +8.31%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.096'
$ws.Range("E43").Value = '  -3.68%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9114'
$ws.Range("E44").Value = '  +3.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4446'
$ws.Range("E45").Value = '  -1.14%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '108.11'
$ws.Range("E46").Value = '  -1.09%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9989'
$ws.Range("E47").Value = '  -0.24%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.722'
$ws.Range("E48").Value = '  -6.65%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '990.91'
$ws.Range("E49").Value = '  +2.76%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1251'
$ws.Range("E50").Value = '  -0.98%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.30'
$ws.Range("E51").Value = '  +1.94%  '
